$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("methods")

# New shared strings must be introduced in this order so the resulting
# sharedStrings.xml table matches the author's original save order:
#   sediment not sieved, roots and rhizomes included, Curtis_et_al_2022, EA

# E2: new sediment_sieved_flag value
$ws.Range("E2").Value = "sediment not sieved"

# D2: new roots_flag value
$ws.Range("D2").Value = "roots and rhizomes included"

# A2: fix capitalization of study id
$ws.Range("A2").Value = "Curtis_et_al_2022"

# U2: fraction carbon method now EA (was "not specified")
$ws.Range("U2").Value = "EA"

# J2: updated dry bulk density sample volume
$ws.Range("J2").Value = 212.37

# Clear out row 3 entirely (J3, M3 notes removed)
$ws.Range("J3").ClearContents()
$ws.Range("M3").ClearContents()

# Update selection / top-left cell to match author's final view
$ws.Range("U3").Select()
